# Generate Report for Handoff
# Replaces the handoff-run identifiers (old uuid "c1dbf75a-...") with the
# new run's identifiers (new uuid "bc815363-...") across the Overview,
# zh-cn and de-de sheets: the source-file hyperlink, the per-language
# handoff-file hyperlink, and that handoff's timestamp.

$wb = $excel.ActiveWorkbook

$newSourceDisplay = "bc815363-6bb6-46ea-a089-fffb3deeab5c.md"
$newSourceUrl = "https://github.com/OpenLocalizationTest/oltest/blob/1ab8e591dcae717f9cdbfc0c1cbf6e9981ea6cfe/e2e/bc815363-6bb6-46ea-a089-fffb3deeab5c.md"

$configDisplay = ".localization-config"
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/1ab8e591dcae717f9cdbfc0c1cbf6e9981ea6cfe/.localization-config"

# ── Overview sheet: only the source-file hyperlink (A2) + config (A3) ──────
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newSourceDisplay

$wsOverview.Range("A2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $newSourceUrl, "", "", $newSourceDisplay) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $configUrl, "", "", $configDisplay) | Out-Null

# ── zh-cn sheet: source file (A2), handoff file + datetime (C2/D2) ─────────
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$newZhHandoffDisplay = "bc815363-6bb6-46ea-a089-fffb3deeab5c.c02f6bf8e16716faa0109ecfd5583ac26674365b.zh-cn.xlf"
$newZhHandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cbac8f4194d84eecfc1ca699ee702b73d38d0b5c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bc815363-6bb6-46ea-a089-fffb3deeab5c.c02f6bf8e16716faa0109ecfd5583ac26674365b.zh-cn.xlf"
$newZhHandoffDatetime = "2016-03-08 06:32:17"

$wsZhCn.Range("A2").Value = $newSourceDisplay
$wsZhCn.Range("C2").Value = $newZhHandoffDisplay
$wsZhCn.Range("D2").Value = $newZhHandoffDatetime

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $newSourceUrl, "", "", $newSourceDisplay) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), $newZhHandoffUrl, "", "", $newZhHandoffDisplay) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $configUrl, "", "", $configDisplay) | Out-Null

# ── de-de sheet: source file (A2), handoff file + datetime (C2/D2) ─────────
$wsDeDe = $wb.Worksheets.Item("de-de")

$newDeHandoffDisplay = "bc815363-6bb6-46ea-a089-fffb3deeab5c.c02f6bf8e16716faa0109ecfd5583ac26674365b.de-de.xlf"
$newDeHandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/365c9264ea992e1f3984f12504e6a5be2c71adcc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bc815363-6bb6-46ea-a089-fffb3deeab5c.c02f6bf8e16716faa0109ecfd5583ac26674365b.de-de.xlf"
$newDeHandoffDatetime = "2016-03-08 06:32:19"

$wsDeDe.Range("A2").Value = $newSourceDisplay
$wsDeDe.Range("C2").Value = $newDeHandoffDisplay
$wsDeDe.Range("D2").Value = $newDeHandoffDatetime

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $newSourceUrl, "", "", $newSourceDisplay) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), $newDeHandoffUrl, "", "", $newDeHandoffDisplay) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $configUrl, "", "", $configDisplay) | Out-Null
